$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-21T07:18:01.126050+00:00"
$ws.Range("K3").Value = "2025-11-21T07:18:01.126073+00:00"
$ws.Range("K4").Value = "2025-11-21T07:18:03.552946+00:00"
$ws.Range("K5").Value = "2025-11-21T07:18:03.552966+00:00"
$ws.Range("K6").Value = "2025-11-21T07:18:05.501753+00:00"
$ws.Range("K7").Value = "2025-11-21T07:18:07.926509+00:00"
$ws.Range("K8").Value = "2025-11-21T07:18:09.893076+00:00"
$ws.Range("K9").Value = "2025-11-21T07:18:09.893094+00:00"
$ws.Range("K10").Value = "2025-11-21T07:18:09.893103+00:00"
$ws.Range("K11").Value = "2025-11-21T07:18:11.775082+00:00"
$ws.Range("K12").Value = "2025-11-21T07:18:14.141960+00:00"
$ws.Range("K13").Value = "2025-11-21T07:18:16.099426+00:00"
$ws.Range("K14").Value = "2025-11-21T07:18:18.046280+00:00"
$ws.Range("K15").Value = "2025-11-21T07:18:18.046309+00:00"
$ws.Range("K16").Value = "2025-11-21T07:18:18.046327+00:00"
$ws.Range("K17").Value = "2025-11-21T07:18:20.441560+00:00"
$ws.Range("K18").Value = "2025-11-21T07:18:25.270869+00:00"
$ws.Range("K19").Value = "2025-11-21T07:18:25.270898+00:00"
$ws.Range("K20").Value = "2025-11-21T07:18:27.732052+00:00"
$ws.Range("K21").Value = "2025-11-21T07:18:27.732084+00:00"
$ws.Range("K22").Value = "2025-11-21T07:18:27.732092+00:00"
$ws.Range("K23").Value = "2025-11-21T07:18:29.683547+00:00"
$ws.Range("K24").Value = "2025-11-21T07:18:29.683574+00:00"
$ws.Range("K25").Value = "2025-11-21T07:18:32.144628+00:00"
$ws.Range("K26").Value = "2025-11-21T07:18:32.144660+00:00"
$ws.Range("K27").Value = "2025-11-21T07:18:32.144680+00:00"
$ws.Range("K28").Value = "2025-11-21T07:18:34.530998+00:00"
$ws.Range("K29").Value = "2025-11-21T07:18:34.531016+00:00"
$ws.Range("K30").Value = "2025-11-21T07:18:36.338509+00:00"
$ws.Range("K31").Value = "2025-11-21T07:18:36.338530+00:00"
$ws.Range("K32").Value = "2025-11-21T07:18:36.338540+00:00"
$ws.Range("K33").Value = "2025-11-21T07:18:36.338548+00:00"
$ws.Range("K34").Value = "2025-11-21T07:18:36.338556+00:00"
$ws.Range("K35").Value = "2025-11-21T07:18:38.292961+00:00"
$ws.Range("K36").Value = "2025-11-21T07:18:38.292994+00:00"
$ws.Range("K37").Value = "2025-11-21T07:18:42.775204+00:00"
$ws.Range("K38").Value = "2025-11-21T07:18:42.775233+00:00"
$ws.Range("K39").Value = "2025-11-21T07:18:42.775251+00:00"
$ws.Range("K40").Value = "2025-11-21T07:18:45.096816+00:00"
